# Remove the "Logging and validation for alfa testing" bullet paragraph
# entirely (whole <w:p>), matching the author's fix described as:
# "Fixed falling by ValueError in float() in value line edits."

$d = $word.ActiveDocument

$targetText = "Logging and validation for alfa testing"

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq $targetText) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Delete the paragraph's whole range, including its paragraph mark,
    # so the entire <w:p> element is removed (not just its text).
    $target.Range.Delete()
}
